$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 665
$ws.Range("J17").Value = 665
$ws.Range("L17").Value = 1995
$ws.Range("N17").Value = -2331
$ws.Range("H70").Value = 3474.875
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 3542.7144
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 10628.1432
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -11168.1432
$ws.Range("H73").Value = 3474.875
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 3542.7144
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 10628.1432
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -12500.1432
$ws.Range("H88").Value = 1644.6154
$ws.Range("J88").Value = 1504.75
$ws.Range("L88").Value = 1504.75
$ws.Range("N88").Value = -2316.75
$ws.Range("H91").Value = 1644.6154
$ws.Range("J91").Value = 1504.75
$ws.Range("L91").Value = 1504.75
$ws.Range("N91").Value = -4312.75
$ws.Range("H95").Value = 90208
$ws.Range("J95").Value = 90208
$ws.Range("L95").Value = 90208
$ws.Range("N95").Value = -95700
$ws.Range("H97").Value = 1260.5
$ws.Range("J97").Value = 1260.5
$ws.Range("L97").Value = 3781.5
$ws.Range("N97").Value = -4773.5
$ws.Range("H100").Value = 1263.5714
$ws.Range("I100").Value = 1309
$ws.Range("J100").Value = 1150
$ws.Range("K100").Value = 1309
$ws.Range("L100").Value = 1150
$ws.Range("M100").Value = -768
$ws.Range("N100").Value = -2232
$ws.Range("H101").Value = 33334034
$ws.Range("I101").Value = 50000550
$ws.Range("K101").Value = 150001650
$ws.Range("M101").Value = -150000028

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 26.5
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = 18
$ws.Range("K5").Value = 35
$ws.Range("L5").Value = 18
$ws.Range("M5").Value = 77
$ws.Range("N5").Value = -242
$ws.Range("H32").Value = 4057781.8
$ws.Range("I32").Value = 3894325.2
$ws.Range("K32").Value = 3894325.2
$ws.Range("M32").Value = -3894038.2
$ws.Range("H97").Value = 850.2308
$ws.Range("I97").Value = 797
$ws.Range("J97").Value = 1027.6666
$ws.Range("K97").Value = 797
$ws.Range("L97").Value = 1027.6666
$ws.Range("M97").Value = -301
$ws.Range("N97").Value = -2019.6666
$ws.Range("H102").Value = 2045.4
$ws.Range("I102").Value = 2045.4
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2045.4
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -423.4000000000001
$ws.Range("N102").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 26.5
$ws.Range("I4").Value = 35
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 35
$ws.Range("L4").Value = 18
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = -248

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 167
$ws.Range("I22").Value = 100.5
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 100.5
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 249.5
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 2257.1667
$ws.Range("I31").Value = 2247.889
$ws.Range("J31").Value = 2285
$ws.Range("K31").Value = 2247.889
$ws.Range("L31").Value = 2285
$ws.Range("M31").Value = -1952.889
$ws.Range("N31").Value = -2875
$ws.Range("H34").Value = 2257.1667
$ws.Range("I34").Value = 2247.889
$ws.Range("J34").Value = 2285
$ws.Range("K34").Value = 2247.889
$ws.Range("L34").Value = 2285
$ws.Range("M34").Value = -2045.889
$ws.Range("N34").Value = -2689
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H60").Value = 6000
$ws.Range("I60").Value = 7000
$ws.Range("J60").Value = 5000
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 5000
$ws.Range("M60").Value = -6489
$ws.Range("N60").Value = -6022
$ws.Range("H74").Value = 28750
$ws.Range("J74").Value = 28750
$ws.Range("L74").Value = 28750
$ws.Range("N74").Value = -30498
$ws.Range("H77").Value = 28750
$ws.Range("J77").Value = 28750
$ws.Range("L77").Value = 86250
$ws.Range("N77").Value = -94986
$ws.Range("H122").Value = 1467
$ws.Range("I122").Value = 1149.1
$ws.Range("J122").Value = 2526.6667
$ws.Range("K122").Value = 3447.3
$ws.Range("L122").Value = 7580.000100000001
$ws.Range("M122").Value = -997.2999999999997
$ws.Range("N122").Value = -12480.0001
$ws.Range("H132").Value = 4674.8
$ws.Range("I132").Value = 4983.1665
$ws.Range("J132").Value = 4212.25
$ws.Range("K132").Value = 14949.4995
$ws.Range("L132").Value = 12636.75
$ws.Range("M132").Value = -12419.4995
$ws.Range("N132").Value = -17696.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 1993
$ws.Range("I120").Value = 1993
$ws.Range("K120").Value = 5979
$ws.Range("M120").Value = -1141
$ws.Range("H139").Value = 6998
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 569.75
$ws.Range("I140").Value = 569.75
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 1709.25
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 3470.75
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 8361.857
$ws.Range("I141").Value = 8361.857
$ws.Range("K141").Value = 25085.571
$ws.Range("M141").Value = -19905.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1666.6666
$ws.Range("I126").Value = 1666.6666
$ws.Range("K126").Value = 4999.9998
$ws.Range("M126").Value = -2529.9998
$ws.Range("H132").Value = 10021.765
$ws.Range("I132").Value = 10059.8
$ws.Range("K132").Value = 30179.4
$ws.Range("M132").Value = -27649.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8480.308000000001
$ws.Range("I7").Value = 8500
$ws.Range("J7").Value = 8471.556
$ws.Range("K7").Value = 8500
$ws.Range("L7").Value = 8471.556
$ws.Range("M7").Value = -8388
$ws.Range("N7").Value = -8695.556
$ws.Range("H46").Value = 2195.5557
$ws.Range("I46").Value = 955
$ws.Range("K46").Value = 955
$ws.Range("M46").Value = -767
$ws.Range("H82").Value = 3488
$ws.Range("J82").Value = 3488
$ws.Range("L82").Value = 3488
$ws.Range("N82").Value = -4210
$ws.Range("H85").Value = 3488
$ws.Range("J85").Value = 3488
$ws.Range("L85").Value = 3488
$ws.Range("N85").Value = -5984
$ws.Range("H101").Value = 34875
$ws.Range("J101").Value = 34875
$ws.Range("L101").Value = 34875
$ws.Range("N101").Value = -41365
$ws.Range("H109").Value = 54000
$ws.Range("J109").Value = 54000
$ws.Range("L109").Value = 54000
$ws.Range("N109").Value = -56774
$ws.Range("H126").Value = 8480.308000000001
$ws.Range("I126").Value = 8500
$ws.Range("J126").Value = 8471.556
$ws.Range("K126").Value = 25500
$ws.Range("L126").Value = 25414.668
$ws.Range("M126").Value = -23030
$ws.Range("N126").Value = -30354.668
$ws.Range("H132").Value = 3298.5386
$ws.Range("J132").Value = 4475.75
$ws.Range("L132").Value = 13427.25
$ws.Range("N132").Value = -18487.25
$ws.Range("H138").Value = 89999.5
$ws.Range("J138").Value = 89999.5
$ws.Range("L138").Value = 89999.5
$ws.Range("N138").Value = -100279.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 18752.625
$ws.Range("I74").Value = 9999
$ws.Range("J74").Value = 20003.143
$ws.Range("K74").Value = 9999
$ws.Range("L74").Value = 20003.143
$ws.Range("M74").Value = -9063
$ws.Range("N74").Value = -21875.143
$ws.Range("H77").Value = 18752.625
$ws.Range("I77").Value = 9999
$ws.Range("J77").Value = 20003.143
$ws.Range("K77").Value = 29997
$ws.Range("L77").Value = 60009.429
$ws.Range("M77").Value = -25317
$ws.Range("N77").Value = -69369.429
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492
$ws.Range("H107").Value = 1772.6957
$ws.Range("I107").Value = 1703.2106
$ws.Range("K107").Value = 5109.6318
$ws.Range("M107").Value = -3189.6318
